$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing row 3, column H (Credito Activo): 1 -> 0
$ws.Range("H3").Value = 0

# Data to append for rows 5-10
# Columns: A=Indice Producto, B=Indice Cantidad, C=Indice Cliente, D=Precio,
#          E=Precio Costo, F=Precio Flete, G=Credito, H=Credito Activo, I=Ganancia
$newRows = @(
    @(4, 1, 2, 198, 179, 3, 1, 0, 16),
    @(3, 1, 1, 245, 230, 3, 1, 0, 12),
    @(3, 5, 1, 245, 230, 3, 1, 1, 60),
    @(4, 4, 2, 198, 179, 3, 1, 0, 64),
    @(1, 2, 2, 198, 190, 3, 1, 0, 10),
    @(4, 7, 2, 198, 179, 3, 1, 1, 0)
)

$startRow = 5
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowData = $newRows[$i]
    for ($col = 1; $col -le $rowData.Count; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowData[$col - 1]
    }
}
